$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Manager")

# Mark a batch of previously "Not Executed" / "FAIL" test cases as executed.
$ws.Range("I13").Value = "PASS"
$ws.Range("I14").Value = "PASS"
$ws.Range("I15").Value = "PASS"
$ws.Range("I16").Value = "PASS"
$ws.Range("I17").Value = "PASS"
$ws.Range("I18").Value = "PASS"
$ws.Range("I19").Value = "PASS"
$ws.Range("I20").Value = "PASS"

# TC_New Customer_04 (row 21) actually failed - record the actual result and status.
$ws.Range("H21").Value = "expected message ""can not start with space"", Actual ""numbers are not allowed"""
$ws.Range("I21").Value = "FAIL"

$ws.Range("I22").Value = "PASS"

# Reflect the reviewer scrolling down to inspect the newly-updated rows.
$ws.Range("I22").Select()
$excel.ActiveWindow.Zoom = 87
